# ---------------------------------------------------------------------------
# Language.xlsx edit: split the single "Sheet1" into 5 sheets (Comm, Property,
# Guild, Tip, Item), add the new localisation strings and populate each
# sheet's ID/English/Chinese table.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# 1. Rename the existing sheet to "Comm" and create the remaining sheets in
#    the right order right after it.
$shComm = $wb.Worksheets.Item(1)
$shComm.Name = "Comm"

$shProperty = $wb.Worksheets.Add($null, $shComm)
$shProperty.Name = "Property"

$shGuild = $wb.Worksheets.Add($null, $shProperty)
$shGuild.Name = "Guild"

$shTip = $wb.Worksheets.Add($null, $shGuild)
$shTip.Name = "Tip"

$shItem = $wb.Worksheets.Add($null, $shTip)
$shItem.Name = "Item"

# ---------------------------------------------------------------------------
# 2. "Comm" sheet (was Sheet1) - header + 6 language rows + 5 blank styled
#    rows reserved for future entries.
# ---------------------------------------------------------------------------
$shComm.Range("A1").Value = "ID"
$shComm.Range("B1").Value = "English"
$shComm.Range("C1").Value = "Chinese"

$shComm.Range("A2").Value = "Langage_Comm_1"
$shComm.Range("B2").Value = "Langage_1"
$shComm.Range("C2").Value = "确认"

$shComm.Range("A3").Value = "Langage_Comm_2"
$shComm.Range("B3").Value = "Langage_2"
$shComm.Range("C3").Value = "取消"

$shComm.Range("A4").Value = "Langage_Comm_3"
$shComm.Range("B4").Value = "Langage_3"
$shComm.Range("C4").Value = "登录"

$shComm.Range("A5").Value = "Langage_Comm_4"
$shComm.Range("B5").Value = "Langage_4"
$shComm.Range("C5").Value = "创建角色"

$shComm.Range("A6").Value = "Langage_Comm_5"
$shComm.Range("B6").Value = "Langage_5"
$shComm.Range("C6").Value = "进入游戏"

$shComm.Range("A7").Value = "Langage_Comm_6"
$shComm.Range("B7").Value = "Langage_6"
$shComm.Range("C7").Value = "中文_6"

# Carry the data-row style down into 5 new blank rows (A8:C12)
$shComm.Range("A2:C2").Copy()
$shComm.Range("A8:C12").PasteSpecial(-4122)
$shComm.Application.CutCopyMode = $false

# Column widths (characters) - chosen so the stored OOXML width matches the
# target (31.875 / 24.5 / 23) as closely as the width->pixel rounding allows.
$shComm.Columns.Item(1).ColumnWidth = 31.142857142857142
$shComm.Columns.Item(2).ColumnWidth = 23.857142857142858
$shComm.Columns.Item(3).ColumnWidth = 22.285714285714285

$shComm.Range("C8").Select()

# ---------------------------------------------------------------------------
# 3. "Property" sheet - single-column (ID only) table of attribute names.
# ---------------------------------------------------------------------------
$shProperty.Range("A1").Value = "ID"
$shProperty.Range("B1").Value = "English"
$shProperty.Range("C1").Value = "Chinese"

$shProperty.Range("A2").Value = "Langage_HP"
$shProperty.Range("A3").Value = "Langage_MAXHP"
$shProperty.Range("A4").Value = "Langage_MP"
$shProperty.Range("A5").Value = "Langage_MAXMP"
$shProperty.Range("A6").Value = "Langage_VP"
$shProperty.Range("A7").Value = "Langage_ATTACK"

$shProperty.Range("A2").Copy()
$shProperty.Range("A8:A28").PasteSpecial(-4122)
$shProperty.Application.CutCopyMode = $false

$shProperty.Columns.Item(1).ColumnWidth = 50.57142857142857

$shProperty.Range("A1:XFD1").Select()

# ---------------------------------------------------------------------------
# 4. "Guild" sheet - header + 1 language row + blank styled rows.
# ---------------------------------------------------------------------------
$shGuild.Range("A1").Value = "ID"
$shGuild.Range("B1").Value = "English"
$shGuild.Range("C1").Value = "Chinese"

$shGuild.Range("A2").Value = "Langage_Guild_1"
$shGuild.Range("B2").Value = "Langage_1"
$shGuild.Range("C2").Value = "确认要加入这个公会吗？点击确认加入"

$shGuild.Range("A2:C2").Copy()
$shGuild.Range("A3:C12").PasteSpecial(-4122)
$shGuild.Range("A16:C16").PasteSpecial(-4122)
$shGuild.Application.CutCopyMode = $false

$shGuild.Columns.Item(1).ColumnWidth = 31.142857142857142
$shGuild.Columns.Item(2).ColumnWidth = 23.857142857142858
$shGuild.Columns.Item(3).ColumnWidth = 22.285714285714285

$shGuild.Range("A22").Select()

# ---------------------------------------------------------------------------
# 5. "Tip" sheet - header row only (placeholder sheet for future tips).
# ---------------------------------------------------------------------------
$shTip.Range("A1").Value = "ID"
$shTip.Range("B1").Value = "English"
$shTip.Range("C1").Value = "Chinese"
$shTip.Range("A1:XFD1").Select()

# ---------------------------------------------------------------------------
# 6. "Item" sheet - header row only (placeholder sheet for future items).
# ---------------------------------------------------------------------------
$shItem.Range("A1").Value = "ID"
$shItem.Range("B1").Value = "English"
$shItem.Range("C1").Value = "Chinese"
$shItem.Range("A1:XFD1").Select()

# Leave the user back on the first sheet, matching the tabSelected="1" state
# recorded for "Comm" in the target workbook.
$shComm.Select()
